$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2024-12-19 00:23:35", -0.1193290994249401, -0.001663952897997995, 0.007942320032144798)
    ,@("2024-12-19 00:23:36", -0.1193290994249401, -0.001647554037813997, 0.007864045583450715)
    ,@("2024-12-19 00:23:37", -0.1185873408327461, -0.001822373985577996, 0.008644419398098709)
    ,@("2024-12-19 00:23:38", -0.1200122981282766, -0.001763510175719995, 0.008465716358430303)
    ,@("2024-12-19 00:23:39", -0.1186621673573973, -0.001687032775293996, 0.008007478620774019)
    ,@("2024-12-19 00:23:40", -0.1185710741969524, -0.001803849347221996, 0.008555374191583337)
    ,@("2024-12-19 00:23:41", -0.1197292586654657, -0.001715629553083995, 0.008216442181412443)
    ,@("2024-12-19 00:23:42", -0.1181025950860931, -0.001595371245067996, 0.007536699366730467)
    ,@("2024-12-19 00:23:43", -0.1183986478575389, -0.001851780583623997, 0.008769932689197022)
    ,@("2024-12-19 00:23:44", -0.1177544890801073, -0.001745997812683996, 0.008223963214703553)
    ,@("2024-12-19 00:23:45", -0.1185515542339999, -0.001770899785555995, 0.008397716878812814)
    ,@("2024-12-19 00:23:46", -0.1193160861163051, -0.001740430298423997, 0.008306453254647364)
    ,@("2024-12-19 00:23:47", -0.1200936313072452, -0.001729295269903996, 0.008307093942608545)
    ,@("2024-12-19 00:23:48", -0.1198854183690855, -0.001732990074821996, 0.008310409605980314)
    ,@("2024-12-19 00:23:49", -0.1212127758498536, -0.001642543274979996, 0.00796388919255339)
    ,@("2024-12-19 00:23:50", -0.119368139350845, -0.001559182402377996, 0.007444668091217676)
    ,@("2024-12-19 00:23:51", -0.1202985909182462, -0.001697256756025995, 0.008167103847056031)
    ,@("2024-12-19 00:23:52", -0.1194136859310675, -0.001716743055935996, 0.008200104644235281)
    ,@("2024-12-19 00:23:53", -0.1193030728076701, -0.001577302130605996, 0.007527079637095211)
)

$startRow = 149
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

Write-Output "Added $($data.Count) rows"